# Apply trade-log update: MarketMaking trade #44 closes (early_exit), and two
# new open trades are appended (momentum #73, EMAArbitrage #74), updating the
# Summary / Strategy Status roll-ups and all affected per-strategy sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value = 1499.7    # Current Capital
$summary.Cells.Item(4, 2).Value = 0.8       # Total P&L $
$summary.Cells.Item(6, 2).Value = 43        # Total Trades
$summary.Cells.Item(7, 2).Value = 24        # Winning Trades
$summary.Cells.Item(9, 2).Value = 55.81     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(6, 3).Value = 99.7       # Capital
$status.Cells.Item(6, 4).Value = 14         # Trades
$status.Cells.Item(6, 5).Value = -0.11      # P&L $
$status.Cells.Item(6, 6).Value = -0.3       # P&L %
$status.Cells.Item(6, 7).Value = 57.14      # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet - close out MarketMaking trade #44 (row 45)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(45, 7).Value = 0.43         # Exit Price
$allTrades.Cells.Item(45, 8).Value = "CLOSED"     # Status
$allTrades.Cells.Item(45, 9).Value = 7.5          # P&L %
$allTrades.Cells.Item(45, 10).Value = 0.03        # P&L $
$allTrades.Cells.Item(45, 11).Value = 99.7        # Capital After
$allTrades.Cells.Item(45, 12).Value = "early_exit" # Exit Reason
$allTrades.Cells.Item(45, 13).Value = 0.14        # Duration (min)

# New row 74: momentum trade #73 (OPEN)
$allTrades.Cells.Item(74, 1).Value = 73
$allTrades.Cells.Item(74, 2).NumberFormat = "@"
$allTrades.Cells.Item(74, 2).Value = "2026-02-18"
$allTrades.Cells.Item(74, 3).Value = "00:10:59"
$allTrades.Cells.Item(74, 4).Value = "momentum"
$allTrades.Cells.Item(74, 5).Value = "UP"
$allTrades.Cells.Item(74, 6).Value = 0.4
$allTrades.Cells.Item(74, 8).Value = "OPEN"
$allTrades.Cells.Item(74, 9).Value = 0
$allTrades.Cells.Item(74, 10).Value = 0
$allTrades.Cells.Item(74, 11).Value = 100
$allTrades.Cells.Item(74, 13).Value = 0
$allTrades.Cells.Item(74, 14).Value = 0
$allTrades.Cells.Item(74, 15).Value = 0
$allTrades.Cells.Item(74, 16).Value = 0.9
$allTrades.Cells.Item(74, 17).Value = "Upward momentum: 21.687% over 10 samples"

# New row 75: EMAArbitrage trade #74 (OPEN)
$allTrades.Cells.Item(75, 1).Value = 74
$allTrades.Cells.Item(75, 2).NumberFormat = "@"
$allTrades.Cells.Item(75, 2).Value = "2026-02-18"
$allTrades.Cells.Item(75, 3).Value = "00:11:00"
$allTrades.Cells.Item(75, 4).Value = "EMAArbitrage"
$allTrades.Cells.Item(75, 5).Value = "UP"
$allTrades.Cells.Item(75, 6).Value = 0.41
$allTrades.Cells.Item(75, 8).Value = "OPEN"
$allTrades.Cells.Item(75, 9).Value = 0
$allTrades.Cells.Item(75, 10).Value = 0
$allTrades.Cells.Item(75, 11).Value = 100
$allTrades.Cells.Item(75, 13).Value = 0
$allTrades.Cells.Item(75, 14).Value = 0
$allTrades.Cells.Item(75, 15).Value = 0
$allTrades.Cells.Item(75, 16).Value = 0.7169
$allTrades.Cells.Item(75, 17).Value = "EMA:up, RSI:50.0, ROC:21.69% | 2/3 UP"

# ---------------------------------------------------------------------------
# momentum sheet - append trade #73 (row 10)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(10, 1).Value = 73
$momentum.Cells.Item(10, 2).NumberFormat = "@"
$momentum.Cells.Item(10, 2).Value = "2026-02-18"
$momentum.Cells.Item(10, 3).Value = "00:10:59"
$momentum.Cells.Item(10, 4).Value = "momentum"
$momentum.Cells.Item(10, 5).Value = "UP"
$momentum.Cells.Item(10, 6).Value = 0.4
$momentum.Cells.Item(10, 8).Value = "OPEN"
$momentum.Cells.Item(10, 9).Value = 0
$momentum.Cells.Item(10, 10).Value = 0
$momentum.Cells.Item(10, 11).Value = 100
$momentum.Cells.Item(10, 12).Value = 0
$momentum.Cells.Item(10, 13).Value = 0
$momentum.Cells.Item(10, 14).Value = 0.9
$momentum.Cells.Item(10, 15).Value = "Upward momentum: 21.687% over 10 samples"
$momentum.Cells.Item(10, 17).Value = 0

# ---------------------------------------------------------------------------
# MarketMaking sheet - close out trade #44 (row 16)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(16, 7).Value = 0.43
$marketMaking.Cells.Item(16, 8).Value = "CLOSED"
$marketMaking.Cells.Item(16, 9).Value = 7.5
$marketMaking.Cells.Item(16, 10).Value = 0.03
$marketMaking.Cells.Item(16, 11).Value = 99.7
$marketMaking.Cells.Item(16, 16).Value = "early_exit"
$marketMaking.Cells.Item(16, 17).Value = 0.14

# ---------------------------------------------------------------------------
# EMAArbitrage sheet - append trade #74 (row 5)
# ---------------------------------------------------------------------------
$emaArb = $wb.Worksheets.Item("EMAArbitrage")
$emaArb.Cells.Item(5, 1).Value = 74
$emaArb.Cells.Item(5, 2).NumberFormat = "@"
$emaArb.Cells.Item(5, 2).Value = "2026-02-18"
$emaArb.Cells.Item(5, 3).Value = "00:11:00"
$emaArb.Cells.Item(5, 4).Value = "EMAArbitrage"
$emaArb.Cells.Item(5, 5).Value = "UP"
$emaArb.Cells.Item(5, 6).Value = 0.41
$emaArb.Cells.Item(5, 8).Value = "OPEN"
$emaArb.Cells.Item(5, 9).Value = 0
$emaArb.Cells.Item(5, 10).Value = 0
$emaArb.Cells.Item(5, 11).Value = 100
$emaArb.Cells.Item(5, 12).Value = 0
$emaArb.Cells.Item(5, 13).Value = 0
$emaArb.Cells.Item(5, 14).Value = 0.7169
$emaArb.Cells.Item(5, 15).Value = "EMA:up, RSI:50.0, ROC:21.69% | 2/3 UP"
$emaArb.Cells.Item(5, 17).Value = 0

Write-Output "edit complete"
